# "exel sheet data test case"
# Adds two new login rows (row 5 = Sufiyan@gmail.com/fourth,
# row 6 = Mumtaz@gmail.com/fifth) below the existing data, each with a
# mailto hyperlink on the column-A email cell (matching the existing
# rows' pattern), then leaves the selection on B6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: Sufiyan@gmail.com / fourth
$ws.Range("A5").Value = "Sufiyan@gmail.com"
$ws.Hyperlinks.Add($ws.Range("A5"), "mailto:Sufiyan@gmail.com")
$ws.Range("A5").Style = "Hyperlink"
$ws.Range("B5").Value = "fourth"

# Row 6: Mumtaz@gmail.com / fifth
$ws.Range("A6").Value = "Mumtaz@gmail.com"
$ws.Hyperlinks.Add($ws.Range("A6"), "mailto:Mumtaz@gmail.com")
$ws.Range("A6").Style = "Hyperlink"
$ws.Range("B6").Value = "fifth"

# Match the saved selection in the committed workbook.
$null = $ws.Range("B6").Select()
